$d = $word.ActiveDocument

# 1. Turn the plain-text dafont URL (end of the "Fonts - " paragraph) into a
#    real hyperlink, matching the style used by every other link in the doc.
$fontsUrl = "https://www.dafont.com/theme.php?cat=101&page=9"
$rng = $d.Content
$found = $rng.Find.Execute($fontsUrl, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $d.Hyperlinks.Add($rng, $fontsUrl)
}

# 2. Add a new paragraph right after it for the pixabay "bones" sound-effects
#    link (kept as plain text, same left indent as the paragraph above it).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newLast = $d.Paragraphs.Last
$newLast.Range.Text = "https://pixabay.com/sound-effects/search/bones/"
